# Apply changes to match target diff for merged_analysis sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B7').Value = 97
$ws.Range('C7').Value = 99
$ws.Range('D7').Value = 96
$ws.Range('E7').Value = 96
$ws.Range('B8').Value = 0
$ws.Range('C8').Value = 14
$ws.Range('D8').Value = 32
$ws.Range('E8').Value = 35
$ws.Range('A9').Value = 'Total distance covered (km)'
$ws.Range('A10').Value = 'Total energy consumption(WH/KM)'
$ws.Range('A11').Value = 'Total SOC consumed(%)'
$ws.Range('A13').Value = 'Peak Power(kW)'
$ws.Range('A14').Value = 'Average Power(kW)'
$ws.Range('A15').Value = 'Total Energy Regenerated(kWh)'
$ws.Range('A16').Value = 'Regenerative Effectiveness(%)'
$ws.Range('B16').Value = 1.201109951715442
$ws.Range('C16').Value = 5.228834716066615
$ws.Range('D16').Value = 4.108522701697714
$ws.Range('E16').Value = 4.266675008452464
$ws.Range('A17').Value = 'Highest Cell Voltage(V)'
$ws.Range('B17').Value = 3.339
$ws.Range('C17').Value = 3.34
$ws.Range('D17').Value = 3.334
$ws.Range('E17').Value = 3.33
$ws.Range('A18').Value = 'Lowest Cell Voltage(V)'
$ws.Range('B18').Value = 2.985
$ws.Range('C18').Value = 3.107
$ws.Range('D18').Value = 3.09
$ws.Range('E18').Value = 3.082
$ws.Range('A19').Value = 'Difference in Cell Voltage(V)'
$ws.Range('A20').Value = 'Minimum Temperature(C)'
$ws.Range('A21').Value = 'Maximum Temperature(C)'
$ws.Range('A22').Value = 'Difference in Temperature(C)'
$ws.Range('B22').Value = 21
$ws.Range('C22').Value = 7
$ws.Range('D22').Value = 10
$ws.Range('E22').Value = 7
$ws.Range('A23').Value = 'Maximum Fet Temperature-BMS(C)'
$ws.Range('A24').Value = 'Maximum Afe Temperature-BMS(C)'
$ws.Range('A25').Value = 'Maximum PCB Temperature-BMS(C)'
$ws.Range('A26').Value = 'Maximum MCU Temperature(C)'
$ws.Range('A27').Value = 'Maximum Motor Temperature(C)'
$ws.Range('A28').Value = 'Abnormal Motor Temperature Detected(C)'
$ws.Range('A29').Value = 'highest cell temp(C)'
$ws.Range('A30').Value = 'lowest cell temp(C)'
$ws.Range('A31').Value = 'Difference between Highest and Lowest Cell Temperature at 100% SOC(C)'
$ws.Range('A32').Value = 'Battery Voltage(V)'
$ws.Range('B32').Value = 55
$ws.Range('C32').Value = 55
$ws.Range('D32').Value = 55
$ws.Range('E32').Value = 55
$ws.Range('A33').Value = 'Total energy charged(kWh)'
$ws.Range('B33').Value = 1.970930148611111
$ws.Range('C33').Value = 1.851597901388889
$ws.Range('D33').Value = 1.386001252777778
$ws.Range('E33').Value = 1.281246388888889
$ws.Range('A34').Value = 'Electricity consumption units(kW)'
$ws.Range('B34').Value = 0.00000006687194293836812
$ws.Range('C34').Value = 0.00000006182627123281674
$ws.Range('D34').Value = 0.0000001082373764390855
$ws.Range('E34').Value = 0.00000009149145878955218
$ws.Range('A35').Value = 'Idling time percentage'
$ws.Range('B35').Value = 18.29305715381665
$ws.Range('C35').Value = 9.657230379211999
$ws.Range('D35').Value = 3.994024208055123
$ws.Range('E35').Value = 5.959138900315371
$ws.Range('A36').Value = 'Time spent in 0-10 km/h'
$ws.Range('B36').Value = 4.760260836210203
$ws.Range('C36').Value = 9.07494454420619
$ws.Range('D36').Value = 23.91536327327053
$ws.Range('E36').Value = 12.43384066913479
$ws.Range('A37').Value = 'Time spent in 10-20 km/h'
$ws.Range('B37').Value = 1.663470144482803
$ws.Range('C37').Value = 10.38475757895849
$ws.Range('D37').Value = 4.378182261654319
$ws.Range('E37').Value = 3.587001234060057
$ws.Range('A38').Value = 'Time spent in 20-30 km/h'
$ws.Range('B38').Value = 3.11724843370413
$ws.Range('C38').Value = 21.83505862469631
$ws.Range('D38').Value = 8.856977346870332
$ws.Range('E38').Value = 7.061565885095297
$ws.Range('A39').Value = 'Time spent in 30-40 km/h'
$ws.Range('B39').Value = 41.70950006393044
$ws.Range('C39').Value = 41.82291116509982
$ws.Range('D39').Value = 18.68044757462118
$ws.Range('E39').Value = 33.23460852872618
$ws.Range('A40').Value = 'Time spent in 40-50 km/h'
$ws.Range('B40').Value = 28.17286791970336
$ws.Range('C40').Value = 6.743160452096757
$ws.Range('D40').Value = 14.79313393701027
$ws.Range('E40').Value = 22.65734265734266
$ws.Range('A41').Value = 'Time spent in 50-60 km/h'
$ws.Range('B41').Value = 0
$ws.Range('C41').Value = 0.05809654589627126
$ws.Range('D41').Value = 19.07984999542669
$ws.Range('E41').Value = 13.45399698340875
$ws.Range('A42').Value = 'Time spent in 60-70 km/h'
$ws.Range('C42').Value = 0.0250871448188444
$ws.Range('D42').Value = 5.884325741638465
$ws.Range('E42').Value = 0.4607157548334019
$ws.Range('A43').Value = 'Time spent in 70-80 km/h'
$ws.Range('C43').Value = 0
$ws.Range('D43').Value = 0
$ws.Range('E43').Value = 0
$ws.Range('A44').Value = 'Time spent in 80-90 km/h'
$ws.Range('B44').Value = 0
$ws.Range('C44').Value = 0
$ws.Range('D44').Value = 0
$ws.Range('E44').Value = 0
